$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the existing H:K data columns (Vin / Vout / dB / Fase)
$ws.Range("H1").Value = "Vin"
$ws.Range("I1").Value = "Vout"
$ws.Range("J1").Value = "dB"
$ws.Range("K1").Value = "Fase"

# Move the active selection to K2 (matches the saved cursor position)
$ws.Range("K2").Select()

# The date-looking number format accidentally applied to I7 should be General
$ws.Range("I7").NumberFormat = "General"
